# Applies the cryptos price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '47.097.60'
$ws.Range("E2").Value = '  +0.44%  '

# Row 3
$ws.Range("D3").Value = '2.477.68'
$ws.Range("E3").Value = '  +0.03%  '

# Row 4
$ws.Range("E4").Value = '  +0.15%  '

# Row 5
$cell = $ws.Range("D5")
$cell.Value = '''319.56'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.33%  '

# Row 6
$cell = $ws.Range("D6")
$cell.Value = '''107.94'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +2.72%  '

# Row 7
$cell = $ws.Range("D7")
$cell.Value = '''0.520'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.34%  '

# Row 8
$cell = $ws.Range("D8")
$cell.Value = '''0.999'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("E9").Value = '  -1.56%  '

# Row 10
$cell = $ws.Range("D10")
$cell.Value = '''38.70'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +6.97%  '

# Row 11
$cell = $ws.Range("D11")
$cell.Value = '''0.0806'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -1.13%  '

# Row 12
$ws.Range("E12").Value = '  +0.35%  '

# Row 13
$cell = $ws.Range("D13")
$cell.Value = '''18.10'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -0.83%  '

# Row 14
$ws.Range("E14").Value = '  -0.35%  '

# Row 15
$ws.Range("D15").Value = '2.847.86'
$ws.Range("E15").Value = '  -0.54%  '

# Row 16
$ws.Range("D16").Value = '2.471.57'
$ws.Range("E16").Value = '  +0.74%  '

# Row 17
$cell = $ws.Range("D17")
$cell.Value = '''0.843'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.26%  '

# Row 18
$ws.Range("D18").Value = '47.029.18'
$ws.Range("E18").Value = '  +0.61%  '

# Row 19
$cell = $ws.Range("D19")
$cell.Value = '''12.65'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.24%  '

# Row 20
$cell = $ws.Range("D20")
$cell.Value = '''6.60'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.81%  '

# Row 21
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0930'
$ws.Range("E21").Value = '  -0.74%  '

# Row 22
$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range("D22")
$cell.Value = '''2.73'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +14.42%  '

# Row 23
$cell = $ws.Range("D23")
$cell.Value = '''70.20'
$cell.Style = "Normal"

# Row 24
$cell = $ws.Range("D24")
$cell.Value = '''244.68'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -2.12%  '

# Row 25
$cell = $ws.Range("D25")
$cell.Value = '''2.55'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.27%  '

# Row 26
$ws.Range("E26").Value = '  -0.14%  '

# Row 27
$cell = $ws.Range("D27")
$cell.Value = '''25.56'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -2.46%  '

# Row 28
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range("D28")
$cell.Value = '''2.28'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +3.22%  '

# Row 29
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range("D29")
$cell.Value = '''10.01'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +1.62%  '

# Row 30
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell = $ws.Range("D30")
$cell.Value = '''34.90'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.64%  '

# Row 31
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range("D31")
$cell.Value = '''0.136'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +0.24%  '

# Row 32
$cell = $ws.Range("D32")
$cell.Value = '''49.35'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.49%  '

# Row 33
$cell = $ws.Range("D33")
$cell.Value = '''19.86'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.95%  '

# Row 34
$ws.Range("E34").Value = '  +0.34%  '

# Row 35
$cell = $ws.Range("D35")
$cell.Value = '''0.0779'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +1.45%  '

# Row 36
$ws.Range("E36").Value = '  +0.29%  '

# Row 37
$ws.Range("E37").Value = '  +2.11%  '

# Row 38
$cell = $ws.Range("D38")
$cell.Value = '''4.63'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.05%  '

# Row 39
$ws.Range("E39").Value = '  -0.76%  '

# Row 40
$ws.Range("E40").Value = '  -0.15%  '

# Row 41
$ws.Range("E41").Value = '  +0.00%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D42")
$cell.Value = '''21.94'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +4.57%  '

# Row 43
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range("D43")
$cell.Value = '''119.24'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -3.28%  '

# Row 44
$cell = $ws.Range("D44")
$cell.Value = '''0.0294'
$cell.Style = "Normal"

# Row 45
$ws.Range("D45").Value = '1.976.62'
$ws.Range("E45").Value = '  +0.00%  '

# Row 46
$cell = $ws.Range("D46")
$cell.Value = '''2.99'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +0.78%  '

# Row 47
$cell = $ws.Range("D47")
$cell.Value = '''2.00'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -4.74%  '

# Row 48
$cell = $ws.Range("D48")
$cell.Value = '''9.05'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.56%  '

# Row 49
$ws.Range("E49").Value = '  -2.30%  '

# Row 50
$cell = $ws.Range("D50")
$cell.Value = '''5.11'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -4.73%  '

# Row 51
$cell = $ws.Range("D51")
$cell.Value = '''57.12'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +5.09%  '
